$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete original row 2 (the "nhl10@gmail.com" row) so rows shift up by one.
$ws.Rows.Item(2).Delete()

# Remove the hyperlink that used to belong to "nhl13@gmail.com" (it stayed on
# A5 because deleting a row does not renumber hyperlink anchors).
foreach ($hl in @($ws.Hyperlinks)) {
    if ($hl.Address -eq "mailto:nhl13@gmail.com") {
        $hl.Delete()
    }
}

# Clear the now-orphaned row 4 contents (used to be "nhl13@gmail.com" / 123456)
# but keep the hyperlink-style formatting on A4.
$ws.Range("A4:B4").ClearContents()

# Update the active cell selection to match the target.
$ws.Range("F10").Select()
